$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings that can look numeric (e.g. "1.00", "0.0000111").
# Force text format on the cells being rewritten so Excel does not silently
# coerce them into numbers (stripping trailing zeros / precision), matching the
# inline-string text cells used in the source workbook.
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = '71.159.86'
$ws.Cells.Item(2, 5).Value = '  +6.68%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '3.686.05'
$ws.Cells.Item(3, 5).Value = '  +19.31%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  +0.04%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '598.67'
$ws.Cells.Item(5, 5).Value = '  +3.82%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '184.10'
$ws.Cells.Item(6, 5).Value = '  +6.74%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '3.683.99'
$ws.Cells.Item(7, 5).Value = '  +19.32%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +0.04%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '0.536'
$ws.Cells.Item(9, 5).Value = '  +4.35%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +8.60%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '6.63'
$ws.Cells.Item(11, 5).Value = '  +3.99%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '0.500'
$ws.Cells.Item(12, 5).Value = '  +6.07%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '39.85'
$ws.Cells.Item(13, 5).Value = '  +11.49%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +6.64%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.307.85'
$ws.Cells.Item(15, 5).Value = '  +19.50%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '3.710.46'
$ws.Cells.Item(16, 5).Value = '  +20.14%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '71.241.25'
$ws.Cells.Item(17, 5).Value = '  +6.88%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +1.83%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '7.50'
$ws.Cells.Item(19, 5).Value = '  +7.70%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '16.98'
$ws.Cells.Item(20, 5).Value = '  -0.06%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '514.95'
$ws.Cells.Item(21, 5).Value = '  +6.56%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '9.24'
$ws.Cells.Item(22, 5).Value = '  +19.03%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '0.744'
$ws.Cells.Item(23, 5).Value = '  +8.17%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '87.40'
$ws.Cells.Item(24, 5).Value = '  +4.94%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '2.43'
$ws.Cells.Item(25, 5).Value = '  +9.49%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '13.44'
$ws.Cells.Item(26, 5).Value = '  +6.26%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '10.99'
$ws.Cells.Item(27, 5).Value = '  +9.64%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '0.999'
$ws.Cells.Item(28, 5).Value = '  -0.03%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +12.05%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '8.10'
$ws.Cells.Item(30, 5).Value = '  +1.56%  '

# Row 31
$ws.Cells.Item(31, 2).Value = 'PEPE'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(31, 4).Value = '0.0000111'
$ws.Cells.Item(31, 5).Value = '  +19.66%  '

# Row 32
$ws.Cells.Item(32, 2).Value = 'EthereumClassic'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(32, 4).Value = '31.74'
$ws.Cells.Item(32, 5).Value = '  +14.00%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).Value = '2.77'
$ws.Cells.Item(33, 5).Value = '  +7.34%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  +4.16%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '1.00'
$ws.Cells.Item(35, 5).Value = '  +0.00%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '6.12'
$ws.Cells.Item(36, 5).Value = '  +9.88%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  +8.58%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +10.81%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '2.17'
$ws.Cells.Item(39, 5).Value = '  +10.96%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '50.91'
$ws.Cells.Item(40, 5).Value = '  +3.89%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  +4.35%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'Arweave'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Cells.Item(42, 4).Value = '45.14'
$ws.Cells.Item(42, 5).Value = '  -6.31%  '

# Row 43
$ws.Cells.Item(43, 2).Value = 'Maker'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(43, 4).Value = '3.177.41'
$ws.Cells.Item(43, 5).Value = '  +14.15%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '8.82'
$ws.Cells.Item(44, 5).Value = '  +6.82%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '410.00'
$ws.Cells.Item(45, 5).Value = '  +11.35%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '2.79'
$ws.Cells.Item(46, 5).Value = '  +7.10%  '

# Row 47
$ws.Cells.Item(47, 5).Value = '  +6.57%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '28.34'
$ws.Cells.Item(48, 5).Value = '  +16.24%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '137.23'
$ws.Cells.Item(49, 5).Value = '  +2.20%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.02%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  +12.64%  '
